# "add default config for npc"
# - Rename NPC id "Player" (row 5) to "Player_0_0"
# - Add a new "Default" NPC row (row 7) that is a copy of the "Enemy" row (row 6)
# - Widen column W (Prefab) to fit the longer values
# - Move selection to A6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 6 ("Enemy") into row 7, but rename the id to "Default"
$ws.Range("A6:Z6").Copy()
$ws.Range("A7:Z7").PasteSpecial()
$ws.Range("A7").Value = "Default"

# Rename existing "Player" NPC entry to "Player_0_0"
$ws.Range("A5").Value = "Player_0_0"

# Widen the "Prefab" column (W) to fit new/longer values
$ws.Columns("W").ColumnWidth = 62.75

# Move the active selection to A6
$ws.Range("A6").Select()
